# Course Content.docx edit:
#  - Add a new bullet "Average linear deviation " right after the
#    "Oscillator  ratio" bullet (same list style/level as its siblings),
#    carrying the lone "_GoBack" bookmark that Word leaves at the last
#    edited spot (it moves off the old "Geometric Mean (GM)" bullet).

$d = $word.ActiveDocument

# Locate the "Oscillator  ratio" paragraph to insert the new bullet after it.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Oscillator\s+ratio") {
        $anchor = $i
        break
    }
}

$anchorPara = $d.Paragraphs.Item($anchor)
$anchorPara.Range.InsertParagraphAfter()

# The freshly-created (empty) paragraph now sits right after the anchor.
$newPara = $d.Paragraphs.Item($anchor + 1)

# Type the bullet text, plus a throwaway trailing marker character "#" so
# that the position just before it is never "end-of-paragraph minus one"
# (a spot where this host mis-resolves a zero-length Bookmarks.Add target
# to a whole, unrelated paragraph). We delete the marker once the
# bookmark has been anchored.
$newPara.Range.InsertBefore("Average linear deviation #")

$newPara = $d.Paragraphs.Item($anchor + 1)
$paraEnd = $newPara.Range.End
$markerPos = $paraEnd - 2

$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()
